$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "27.660.67"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.859.74"
$ws.Range("E3").Value = "  -1.09%  "
Set-TextValue "D4" "1.014"
$ws.Range("E4").Value = "  +1.00%  "
Set-TextValue "D5" "332.88"
$ws.Range("E5").Value = "  +0.11%  "
Set-TextValue "D6" "1.012"
$ws.Range("E6").Value = "  +0.83%  "
Set-TextValue "D7" "0.4631"
$ws.Range("E7").Value = "  -2.04%  "
Set-TextValue "D8" "0.3876"
$ws.Range("E8").Value = "  -2.09%  "
Set-TextValue "D9" "46.03"
$ws.Range("E9").Value = "  -3.80%  "
Set-TextValue "D10" "0.07944"
$ws.Range("E10").Value = "  -1.54%  "
Set-TextValue "D11" "0.9944"
$ws.Range("E11").Value = "  -3.73%  "
Set-TextValue "D12" "21.44"
$ws.Range("E12").Value = "  -3.61%  "
$ws.Range("D13").Value = "1.870.33"
$ws.Range("E13").Value = "  -0.19%  "
Set-TextValue "D14" "5.959"
$ws.Range("E14").Value = "  -0.39%  "
Set-TextValue "D15" "7.151"
$ws.Range("E15").Value = "  +0.11%  "
Set-TextValue "D16" "1.015"
$ws.Range("E16").Value = "  +0.76%  "
Set-TextValue "D17" "87.68"
$ws.Range("E17").Value = "  +0.43%  "
Set-TextValue "D18" "0.06706"
$ws.Range("E18").Value = "  +0.48%  "
Set-TextValue "D19" "0.00001040"
$ws.Range("E19").Value = "  -1.02%  "
Set-TextValue "D20" "16.85"
$ws.Range("E20").Value = "  -2.74%  "
Set-TextValue "D21" "1.011"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").Value = "27.662.40"
$ws.Range("E22").Value = "  -0.50%  "
Set-TextValue "D23" "5.435"
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("E24").Value = "  -1.46%  "
Set-TextValue "D25" "2.321"
$ws.Range("E25").Value = "  +0.70%  "
Set-TextValue "D26" "158.39"
$ws.Range("E26").Value = "  -0.73%  "
Set-TextValue "D27" "19.65"
$ws.Range("E27").Value = "  -3.10%  "
Set-TextValue "D28" "2.110"
$ws.Range("E28").Value = "  -0.05%  "
Set-TextValue "D29" "5.342"
$ws.Range("E29").Value = "  -4.91%  "
Set-TextValue "D30" "121.16"
$ws.Range("E30").Value = "  -0.95%  "
Set-TextValue "D31" "0.9669"
$ws.Range("E31").Value = "  -1.97%  "
Set-TextValue "D32" "0.09421"
$ws.Range("E32").Value = "  -1.37%  "
Set-TextValue "D33" "3.643"
$ws.Range("E33").Value = "  +1.25%  "
Set-TextValue "D34" "5.277"
$ws.Range("E34").Value = "  -1.95%  "
Set-TextValue "D35" "1.320"
$ws.Range("E35").Value = "  -9.14%  "
Set-TextValue "D36" "0.05993"
$ws.Range("E36").Value = "  -2.24%  "
Set-TextValue "D37" "0.02212"
$ws.Range("E37").Value = "  -2.14%  "
Set-TextValue "D38" "1.196"
$ws.Range("E38").Value = "  -3.50%  "
$ws.Range("B39").Value = "Frax"
$ws.Range("C39").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D39" "1.011"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D40" "8.108"
$ws.Range("E40").Value = "  -0.68%  "
Set-TextValue "D41" "0.5882"
$ws.Range("E41").Value = "  -2.65%  "
Set-TextValue "D42" "0.1873"
$ws.Range("E42").Value = "  -1.76%  "
Set-TextValue "D43" "10.20"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("E44").Value = "  -0.68%  "
Set-TextValue "D45" "0.5595"
$ws.Range("E45").Value = "  -2.69%  "
Set-TextValue "D46" "12.10"
$ws.Range("E46").Value = "  -1.03%  "
Set-TextValue "D47" "1.911"
$ws.Range("E47").Value = "  -2.15%  "
Set-TextValue "D48" "3.285"
$ws.Range("E48").Value = "  -2.84%  "
Set-TextValue "D49" "0.06764"
$ws.Range("E49").Value = "  -2.19%  "
Set-TextValue "D50" "111.91"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D51" "0.00000000298"
$ws.Range("E51").Value = "  -2.91%  "
